$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the forgotten card, pushing the total row down
$ws.Rows("12:12").Insert()

$ws.Range("A12").Value = "Scion of Draco"
$ws.Range("B12").Value = "Modern Horizons 2"
$ws.Range("C12").Value = "Normal"
$ws.Range("D12").Value = 1.99

# Match the number format used by the rest of the Price column
$fmt = '_-* #,##0.00\ [$€-1]_-;\-* #,##0.00\ [$€-1]_-;_-* "-"??\ [$€-1]_-;_-@_-'
$ws.Range("D12").NumberFormat = $fmt

# Extend the total formula to include the new row
$ws.Range("D14").Formula = "=SUM(D2:D12)"

$ws.Range("D13").Select()
